# "Looking at my version of the line function"
# Extend the worked example on the "Alt1" sheet: add blank (but styled)
# cells in column F for the existing G22:I37 block, and extend the
# C:E "In/Out" sample table down from row 34 through row 49.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt1")   # "Alt1" is already the active/selected tab
$ws.Activate() | Out-Null

# --- Column F: add empty, centre-styled cells alongside rows 22-33 ---
# (matches the style already used on F3:F11) by copying formats from F3.
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F22:F33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Columns C:E rows 34-49: extend the In/Out worked-example table ---
$data = @(
    @("In",  "Out", 1),
    @("Out", "In",  1),
    @("In",  "Out", 2),
    @("Out", "In",  2),
    @("In",  "Out", 3),
    @("Out", "In",  3),
    @("In",  "Out", 4),
    @("Out", "In",  4),
    @("In",  "Out", 5),
    @("Out", "In",  6),
    @("In",  "Out", 6),
    @("Out", "Out", 7),
    @("Out", "In",  8),
    @("In",  "In",  9),
    @("In",  "In",  10)
)

$row = 34
foreach ($entry in $data) {
    $ws.Cells.Item($row, 3).Value = $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]
    $ws.Cells.Item($row, 5).Value = $entry[2]
    $row++
}

# Row 49 only has a single "In" value in column C.
$ws.Cells.Item(49, 3).Value = "In"

# --- Sheet view bookkeeping: selection moved while scrolled down ---
$ws.Range("G42").Select() | Out-Null
